$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet: row for b.md (row 3) moves from
# "Handed back: in sync with en-US" / 2016-03-19 03:32:36
# to "Ready for handoff" / 2016-03-19 03:34:05
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-19 03:34:05"

# ---------------------------------------------------------------
# zh-cn sheet: row 3 (b.md) gets a new handoff status/file/datetime
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-19 03:33:57"

foreach ($h in $wsZhCn.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$3') {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------
# de-de sheet: row 3 (b.md) gets a new handoff status/file/datetime
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-19 03:34:05"

foreach ($h in $wsDeDe.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$3') {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
